# Update the two-digit-divided-by-one-digit practice table.
# Each formula cell is addressed by its (row, column) position in the
# single table so that duplicate/overlapping values (e.g. "13÷9=" and
# "95÷7=" occur both as an old value in one cell and a new value in
# another) are never confused with a global text find/replace.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "43÷7="
$t.Cell(1, 2).Range.Text = "38÷2="
$t.Cell(1, 3).Range.Text = "13÷9="
$t.Cell(1, 4).Range.Text = "53÷8="
$t.Cell(1, 5).Range.Text = "64÷8="

# Row 5
$t.Cell(5, 1).Range.Text = "59÷7="
$t.Cell(5, 2).Range.Text = "28÷9="
$t.Cell(5, 3).Range.Text = "65÷9="
$t.Cell(5, 4).Range.Text = "20÷7="
$t.Cell(5, 5).Range.Text = "11÷6="

# Row 9
$t.Cell(9, 1).Range.Text = "63÷5="
$t.Cell(9, 2).Range.Text = "65÷3="
$t.Cell(9, 3).Range.Text = "38÷9="
$t.Cell(9, 4).Range.Text = "58÷7="
$t.Cell(9, 5).Range.Text = "68÷2="

# Row 13
$t.Cell(13, 1).Range.Text = "95÷7="
$t.Cell(13, 2).Range.Text = "25÷5="
$t.Cell(13, 3).Range.Text = "94÷8="
$t.Cell(13, 4).Range.Text = "81÷8="
$t.Cell(13, 5).Range.Text = "76÷9="

# Row 17
$t.Cell(17, 1).Range.Text = "41÷4="
$t.Cell(17, 2).Range.Text = "24÷7="
$t.Cell(17, 3).Range.Text = "30÷7="
$t.Cell(17, 4).Range.Text = "52÷2="
$t.Cell(17, 5).Range.Text = "23÷2="

Write-Output "Updated 25 division problems"
